# Adds a new column D ("IdxPar") to the "Groupes" sheet, filled with the
# constant index 9040000000 for every data row, and marks that whole
# column as unlocked (applyProtection + locked=0), matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new column.
$ws.Range("D1").Value2 = "IdxPar"

# Data rows 2..171 all receive the same constant value.
$ws.Range("D2:D171").Value2 = 9040000000

# New column is unlocked (adds the applyProtection/locked=0 cell style),
# applied to the whole new column range including the header.
$ws.Range("D1:D171").Locked = $False

# New column width.
$ws.Columns.Item(4).ColumnWidth = 18

# Move the active selection to E169 (just right of the new column),
# matching the post-edit workbook's cursor position.
$ws.Range("E169").Select() | Out-Null
